$wb = $excel.ActiveWorkbook

$mainWs = $wb.Worksheets.Item("main")
$modelWs = $wb.Worksheets.Item("model")

$mainWs.Range("O2").Value = 155
$modelWs.Range("I6").Formula = "=H6*1.13"
$modelWs.Range("J6").Formula = "=I6*1.12"
